# Apply the authored edit:
#  - Insert a new slide at position 3: "Comparisons of the first two slides"
#  - Append a new slide at the end (position 6): "Observation of the last two slides"
# Both new slides use the "Titre et contenu" (Title and Content) layout, the
# same layout already used by the two pre-existing content slides.

$p = $ppt.ActivePresentation

# Grab the "Title and Content" custom layout by reading it off an existing
# slide that already uses it (Master.CustomLayouts.Item(N) does not reliably
# index by position in this host, but Slide.CustomLayout correctly resolves
# to the layout actually used by that slide).
$layout = $p.Slides.Item(3).CustomLayout

# --- New slide inserted at position 3 -------------------------------------
$s3 = $p.Slides.AddSlide(3, $layout)

$title3 = $s3.Shapes.Item(1)
$title3.TextFrame.TextRange.Text = "Comparisons of the first two slides"

$body3 = $s3.Shapes.Item(2)
$body3.TextFrame.TextRange.Text = "Looking at the first two slides, we can observe that we get the same sequence between the simulation and the signal analysis."

# --- New slide appended at the end (position 6) ----------------------------
$s6 = $p.Slides.AddSlide($p.Slides.Count + 1, $layout)

$title6 = $s6.Shapes.Item(1)
$title6.TextFrame.TextRange.Text = "Observation of the last two slides"

$body6 = $s6.Shapes.Item(2)
$body6.TextFrame.TextRange.Text = "Looking at the last two slides, we can see that we are reading FF because the sensor is not connected (and 00 when it is connected because we have not set a read delay). We can also see the writing and reading props E0 and E1. Finally, we can see that we are interacting with register 02. "

Write-Host ("Slides.Count = {0}" -f $p.Slides.Count)
